$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as text, preserving its original style/formatting.
# This prevents Excel from auto-converting numeric-looking strings (e.g. "303.89")
# into actual numbers, which would change the cell type and lose the exact text.
function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '23.764.39'
Set-TextValue $ws.Range("E2") '  +1.57%  '

Set-TextValue $ws.Range("D3") '1.657.76'
Set-TextValue $ws.Range("E3") '  +1.69%  '

Set-TextValue $ws.Range("D4") '1.002'
Set-TextValue $ws.Range("E4") '  +0.20%  '

Set-TextValue $ws.Range("E5") '  +0.20%  '

Set-TextValue $ws.Range("D6") '303.89'
Set-TextValue $ws.Range("E6") '  +0.52%  '

Set-TextValue $ws.Range("D7") '0.3806'
Set-TextValue $ws.Range("E7") '  +0.98%  '

Set-TextValue $ws.Range("D8") '0.3632'
Set-TextValue $ws.Range("E8") '  +0.07%  '

Set-TextValue $ws.Range("D9") '51.04'
Set-TextValue $ws.Range("E9") '  -1.74%  '

Set-TextValue $ws.Range("D10") '1.254'
Set-TextValue $ws.Range("E10") '  +2.69%  '

Set-TextValue $ws.Range("D11") '0.08235'
Set-TextValue $ws.Range("E11") '  +0.83%  '

Set-TextValue $ws.Range("D12") '1.003'
Set-TextValue $ws.Range("E12") '  +0.33%  '

Set-TextValue $ws.Range("D13") '22.68'
Set-TextValue $ws.Range("E13") '  +2.08%  '

Set-TextValue $ws.Range("D14") '6.547'
Set-TextValue $ws.Range("E14") '  +1.21%  '

Set-TextValue $ws.Range("D15") '7.478'
Set-TextValue $ws.Range("E15") '  +2.09%  '

Set-TextValue $ws.Range("D16") '0.00001244'
Set-TextValue $ws.Range("E16") '  +0.34%  '

Set-TextValue $ws.Range("D17") '1.657.27'
Set-TextValue $ws.Range("E17") '  +2.30%  '

Set-TextValue $ws.Range("D18") '97.85'
Set-TextValue $ws.Range("E18") '  +3.17%  '

Set-TextValue $ws.Range("D19") '0.06993'
Set-TextValue $ws.Range("E19") '  +0.69%  '

Set-TextValue $ws.Range("D20") '6.815'
Set-TextValue $ws.Range("E20") '  +4.23%  '

Set-TextValue $ws.Range("D21") '17.78'
Set-TextValue $ws.Range("E21") '  +1.35%  '

Set-TextValue $ws.Range("E22") '  +0.12%  '

Set-TextValue $ws.Range("E23") '  +2.56%  '

Set-TextValue $ws.Range("D24") '23.761.20'
Set-TextValue $ws.Range("E24") '  +1.58%  '

Set-TextValue $ws.Range("D25") '2.550'
Set-TextValue $ws.Range("E25") '  +1.58%  '

Set-TextValue $ws.Range("D26") '3.070'
Set-TextValue $ws.Range("E26") '  -0.23%  '

Set-TextValue $ws.Range("E27") '  +0.98%  '

Set-TextValue $ws.Range("D28") '151.68'
Set-TextValue $ws.Range("E28") '  +0.62%  '

Set-TextValue $ws.Range("D29") '5.232'
Set-TextValue $ws.Range("E29") '  -0.80%  '

Set-TextValue $ws.Range("D30") '134.40'
Set-TextValue $ws.Range("E30") '  +1.21%  '

Set-TextValue $ws.Range("D31") '1.840.42'
Set-TextValue $ws.Range("E31") '  +2.07%  '

Set-TextValue $ws.Range("D32") '6.944'
Set-TextValue $ws.Range("E32") '  +5.00%  '

Set-TextValue $ws.Range("D33") '2.190'
Set-TextValue $ws.Range("E33") '  +2.43%  '

Set-TextValue $ws.Range("E34") '  +1.89%  '

Set-TextValue $ws.Range("D35") '11.83'
Set-TextValue $ws.Range("E35") '  +4.93%  '

Set-TextValue $ws.Range("D36") '0.02824'
Set-TextValue $ws.Range("E36") '  +2.18%  '

Set-TextValue $ws.Range("D37") '0.2527'
Set-TextValue $ws.Range("E37") '  +1.54%  '

Set-TextValue $ws.Range("D38") '6.140'
Set-TextValue $ws.Range("E38") '  +2.93%  '

Set-TextValue $ws.Range("D39") '0.08830'
Set-TextValue $ws.Range("E39") '  +0.84%  '

Set-TextValue $ws.Range("D40") '0.07097'
Set-TextValue $ws.Range("E40") '  -0.49%  '

Set-TextValue $ws.Range("D41") '13.34'
Set-TextValue $ws.Range("E41") '  +11.19%  '

Set-TextValue $ws.Range("D42") '0.7077'
Set-TextValue $ws.Range("E42") '  +1.29%  '

Set-TextValue $ws.Range("E43") '  +0.92%  '

Set-TextValue $ws.Range("D44") '15.91'
Set-TextValue $ws.Range("E44") '  +0.63%  '

Set-TextValue $ws.Range("D45") '0.6553'
Set-TextValue $ws.Range("E45") '  +1.46%  '

Set-TextValue $ws.Range("D46") '2.333'
Set-TextValue $ws.Range("E46") '  +2.57%  '

Set-TextValue $ws.Range("E47") '  +0.23%  '

Set-TextValue $ws.Range("D48") '3.964'
Set-TextValue $ws.Range("E48") '  +0.10%  '

Set-TextValue $ws.Range("D50") '128.21'
Set-TextValue $ws.Range("E50") '  +1.61%  '

Set-TextValue $ws.Range("E51") '  +0.63%  '
